$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 78, shifting existing rows 78..113 down to 79..114
$ws.Range("A78:R78").Insert()

# Populate the newly inserted row 78 with its data
$ws.Range("A78").Value = 4
$ws.Range("B78").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C78").Value = "Los Lagos"
$ws.Range("D78").Value = 44466
$ws.Range("E78").Value = 10
$ws.Range("F78").Value = 100112039
$ws.Range("G78").Value = "Ciboulette"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 90
$ws.Range("K78").Value = 4000
$ws.Range("L78").Value = 4000
$ws.Range("M78").Value = 4000
$ws.Range("N78").Value = "$/docena de atados"
$ws.Range("O78").Value = "Región Metropolitana"
$ws.Range("P78").Value = 1333
$ws.Range("Q78").Value = 3
$ws.Range("R78").Value = "Hortaliza"
